$wb = $excel.ActiveWorkbook

# The CodeSystem "Metadata" sheet holds the generated implementation-guide
# summary table (Property / Value pairs). Re-deploying the guide bumps the
# resource's Status and Date metadata.
$ws = $wb.Worksheets.Item("Metadata")

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Date: 2023-05-12T12:33:13+00:00 -> 2023-08-01T16:12:28+00:00
$ws.Range("B8").Value = "2023-08-01T16:12:28+00:00"
